# -----------------------------------------------------------------------
# Weekly Fruta/Hortaliza refresh for "Feria Lagunitas de Puerto Montt - Ciruela".
# Two brand-new price observations are inserted as rows 71-72 (Black Amber,
# 2022-01-13), which pushes every existing observation that used to live in
# rows 71-115 down by two rows (now rows 73-117).
# -----------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DataRow {
    param($ws, $rowNum, $values)
    $arr = New-Object "object[,]" 1,20
    for ($i = 0; $i -lt 20; $i++) {
        $arr[0,$i] = $values[$i]
    }
    $ws.Range("A" + $rowNum + ":T" + $rowNum).Value = $arr
}

Set-DataRow $ws 71 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44574, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Primera', 200, 18000, 18500, 18250, '$/caja 15 kilos granel', 'Región de O''Higgins', 1217, 15)
Set-DataRow $ws 72 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44574, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Segunda', 100, 16000, 16000, 16000, '$/caja 15 kilos granel', 'Región de O''Higgins', 1067, 15)
Set-DataRow $ws 73 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44299, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Angeleno', 'Primera', 200, 11000, 12000, 11500, '$/caja 14 kilos granel', 'Región de O''Higgins', 821, 14)
Set-DataRow $ws 74 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44544, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Red Beaut', 'Primera', 700, 19000, 20000, 19500, '$/caja 15 kilos granel', 'Región Metropolitana', 1300, 15)
Set-DataRow $ws 75 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44218, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Primera', 200, 13000, 13000, 13000, '$/caja 14 kilos granel', 'Región de O''Higgins', 929, 14)
Set-DataRow $ws 76 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44218, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Segunda', 200, 10000, 10000, 10000, '$/caja 14 kilos granel', 'Región de O''Higgins', 714, 14)
Set-DataRow $ws 77 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44320, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Angeleno', 'Primera', 300, 12000, 12000, 12000, '$/caja 14 kilos granel', 'Región de O''Higgins', 857, 14)
Set-DataRow $ws 78 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44306, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Angeleno', 'Primera', 200, 12000, 13000, 12500, '$/caja 14 kilos granel', 'Región de O''Higgins', 893, 14)
Set-DataRow $ws 79 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44295, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Angeleno', 'Primera', 200, 11000, 12000, 11500, '$/caja 14 kilos granel', 'Región de O''Higgins', 821, 14)
Set-DataRow $ws 80 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44230, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Primera', 100, 13000, 13000, 13000, '$/caja 14 kilos granel', 'Región de O''Higgins', 929, 14)
Set-DataRow $ws 81 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44230, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Segunda', 100, 10000, 10000, 10000, '$/caja 14 kilos granel', 'Región de O''Higgins', 714, 14)
Set-DataRow $ws 82 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44316, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Angeleno', 'Primera', 300, 11000, 12000, 11500, '$/caja 14 kilos granel', 'Región de O''Higgins', 821, 14)
Set-DataRow $ws 83 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44211, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Primera', 150, 13000, 13000, 13000, '$/caja 14 kilos granel', 'Región de O''Higgins', 929, 14)
Set-DataRow $ws 84 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44211, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Segunda', 150, 10000, 10000, 10000, '$/caja 14 kilos granel', 'Región de O''Higgins', 714, 14)
Set-DataRow $ws 85 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44211, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Lemon', 'Primera', 400, 14000, 14500, 14250, '$/caja 14 kilos granel', 'Región de O''Higgins', 1018, 14)
Set-DataRow $ws 86 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44313, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Angeleno', 'Primera', 300, 11000, 12000, 11500, '$/caja 14 kilos granel', 'Región de O''Higgins', 821, 14)
Set-DataRow $ws 87 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44334, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Angeleno', 'Primera', 300, 11000, 12000, 11500, '$/caja 14 kilos granel', 'Región de O''Higgins', 821, 14)
Set-DataRow $ws 88 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44330, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Angeleno', 'Primera', 300, 11000, 12000, 11500, '$/caja 14 kilos granel', 'Región de O''Higgins', 821, 14)
Set-DataRow $ws 89 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44196, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Red Beaut', 'Primera', 100, 15000, 15000, 15000, '$/caja 14 kilos granel', 'Región Metropolitana', 1071, 14)
Set-DataRow $ws 90 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44196, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Red Beaut', 'Segunda', 100, 12000, 12000, 12000, '$/caja 14 kilos granel', 'Región Metropolitana', 857, 14)
Set-DataRow $ws 91 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44239, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Primera', 300, 12000, 13000, 12500, '$/caja 14 kilos granel', 'Región de O''Higgins', 893, 14)
Set-DataRow $ws 92 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44239, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Segunda', 150, 9000, 9000, 9000, '$/caja 14 kilos granel', 'Región de O''Higgins', 643, 14)
Set-DataRow $ws 93 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44208, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Primera', 150, 13000, 13000, 13000, '$/caja 14 kilos granel', 'Región de O''Higgins', 929, 14)
Set-DataRow $ws 94 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44208, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Segunda', 150, 10000, 10000, 10000, '$/caja 14 kilos granel', 'Región de O''Higgins', 714, 14)
Set-DataRow $ws 95 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44250, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Primera', 200, 11000, 12000, 11500, '$/caja 14 kilos granel', 'Región de O''Higgins', 821, 14)
Set-DataRow $ws 96 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44250, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Segunda', 150, 8000, 8000, 8000, '$/caja 14 kilos granel', 'Región de O''Higgins', 571, 14)
Set-DataRow $ws 97 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44560, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Red Beaut', 'Primera', 400, 18000, 19000, 18500, '$/caja 15 kilos granel', 'Región Metropolitana', 1233, 15)
Set-DataRow $ws 98 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44560, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Red Beaut', 'Segunda', 200, 16000, 16000, 16000, '$/caja 15 kilos granel', 'Región Metropolitana', 1067, 15)
Set-DataRow $ws 99 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44229, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Primera', 150, 13000, 13000, 13000, '$/caja 14 kilos granel', 'Región de O''Higgins', 929, 14)
Set-DataRow $ws 100 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44229, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Segunda', 150, 10000, 10000, 10000, '$/caja 14 kilos granel', 'Región de O''Higgins', 714, 14)
Set-DataRow $ws 101 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44298, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Angeleno', 'Primera', 100, 11000, 12000, 11500, '$/caja 14 kilos granel', 'Región de O''Higgins', 821, 14)
Set-DataRow $ws 102 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44558, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Red Beaut', 'Primera', 500, 18000, 19000, 18500, '$/caja 15 kilos granel', 'Región Metropolitana', 1233, 15)
Set-DataRow $ws 103 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44558, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Red Beaut', 'Segunda', 250, 16000, 16000, 16000, '$/caja 15 kilos granel', 'Región Metropolitana', 1067, 15)
Set-DataRow $ws 104 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44301, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Angeleno', 'Primera', 200, 12000, 13000, 12500, '$/caja 14 kilos granel', 'Región de O''Higgins', 893, 14)
Set-DataRow $ws 105 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44278, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Primera', 300, 11000, 12000, 11500, '$/caja 14 kilos granel', 'Región de O''Higgins', 821, 14)
Set-DataRow $ws 106 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44322, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Angeleno', 'Primera', 200, 11000, 12000, 11500, '$/caja 14 kilos granel', 'Región de O''Higgins', 821, 14)
Set-DataRow $ws 107 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44194, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Red Beaut', 'Primera', 200, 15000, 15000, 15000, '$/caja 14 kilos granel', 'Región Metropolitana', 1071, 14)
Set-DataRow $ws 108 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44194, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Red Beaut', 'Segunda', 200, 12000, 12000, 12000, '$/caja 14 kilos granel', 'Región Metropolitana', 857, 14)
Set-DataRow $ws 109 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44236, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Primera', 500, 13000, 14000, 13500, '$/caja 14 kilos granel', 'Región de O''Higgins', 964, 14)
Set-DataRow $ws 110 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44236, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Segunda', 100, 10000, 10000, 10000, '$/caja 14 kilos granel', 'Región de O''Higgins', 714, 14)
Set-DataRow $ws 111 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44293, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Angeleno', 'Primera', 200, 11000, 12000, 11500, '$/caja 14 kilos granel', 'Región de O''Higgins', 821, 14)
Set-DataRow $ws 112 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44266, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Primera', 200, 11000, 12000, 11500, '$/caja 14 kilos granel', 'Región de O''Higgins', 821, 14)
Set-DataRow $ws 113 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44266, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Segunda', 100, 8000, 8000, 8000, '$/caja 14 kilos granel', 'Región de O''Higgins', 571, 14)
Set-DataRow $ws 114 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44533, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Red Beaut', 'Primera', 600, 22000, 23000, 22500, '$/caja 15 kilos granel', 'Región Metropolitana', 1500, 15)
Set-DataRow $ws 115 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44264, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Primera', 200, 11000, 12000, 11500, '$/caja 14 kilos granel', 'Región de O''Higgins', 821, 14)
Set-DataRow $ws 116 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44264, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Black Amber', 'Segunda', 100, 8000, 8000, 8000, '$/caja 14 kilos granel', 'Región de O''Higgins', 571, 14)
Set-DataRow $ws 117 @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44327, 10, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103002, 'Ciruela', 'Angeleno', 'Primera', 300, 11000, 12000, 11500, '$/caja 14 kilos granel', 'Región de O''Higgins', 821, 14)

# New rows (116, 117) do not inherit the "yyyy-mm-dd" style used by column D,
# so make sure the date cells keep the same number format as the rest of the
# column.
$ws.Range("D116").NumberFormat = $ws.Range("D115").NumberFormat()
$ws.Range("D117").NumberFormat = $ws.Range("D115").NumberFormat()
